# november.2025.xlsx -- "Add files via upload"
#
# What changed (per the OOXML diff):
#   1. "Support" sheet: the saved view scrolled so column D is the left-most
#      visible column, and the selection became the header row range A1:L1
#      (previously the lingering selection was the single cell K9).
#   2. "Delivery" sheet: the saved selection moved from E3 to E6.
#   3. "Delivery" sheet, row 1 (the merged two-row-per-person header): the
#      "Support" / "Shift Timings" header cells were pointing at a stray,
#      slightly-mistyped duplicate pair of shared strings ("Support" with no
#      trailing space / "Shitt Timings"). They get repointed at the correct,
#      already-used pair of strings ("Support " / "Shift Timings") -- the
#      same ones the header row on the "Support" sheet already used -- and
#      cell E1 loses its bold weight to match the rest of that header run.
#      Once nothing references the stray duplicate strings any more, they
#      drop out of the shared-string table on save.

$wb = $excel.ActiveWorkbook

$support  = $wb.Worksheets.Item("Support")
$delivery = $wb.Worksheets.Item("Delivery")

# --- 1. Support sheet: scroll + selection -----------------------------
$support.Activate()
$support.Range("A1:L1").Select()
$excel.ActiveWindow.ScrollColumn = 4   # column D becomes the left edge
$excel.ActiveWindow.ScrollRow = 1

# --- 2 & 3. Delivery sheet: selection + header row fix-up --------------
$delivery.Activate()

$delivery.Range("C1").Value = "Support "
$delivery.Range("D1").Value = "Shift Timings"
$delivery.Range("E1").Value = "Support "
$delivery.Range("E1").Font.Bold = $false
$delivery.Range("F1").Value = "Shift Timings"
$delivery.Range("G1").Value = "Support "
$delivery.Range("H1").Value = "Shift Timings"
$delivery.Range("I1").Value = "Support "
$delivery.Range("J1").Value = "Shift Timings"
$delivery.Range("K1").Value = "Support "
$delivery.Range("L1").Value = "Shift Timings"

$delivery.Range("E6").Select()
